# Update exported Selenium and unit test case reports
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data set for rows 1..31 (A: Filename, B: Test Case Name, C: Line Number)
$data = @(
    @("Filename", "Test Case Name", 0),
    @("inventory_test.py", "test_add_items", 4),
    @("inventory_test.py", "test_edit_item_quantity", 11),
    @("inventory_test.py", "test_save_to_file", 17),
    @("test_ingredients_scenarios.py", "test_additional_ingredient_scenarios", 10),
    @("test_ingredients_scenarios.py", "test_chicken_and_rice", 18),
    @("test_ingredients_scenarios.py", "test_chicken_lemon_parsley", 24),
    @("test_ingredients_scenarios.py", "test_empty_ingredient_list", 30),
    @("test_ingredients_scenarios.py", "test_conflicting_restrictions", 36),
    @("test_ingredients_scenarios.py", "test_gibberish_input", 42),
    @("test_bias_logic.py", "test_bias_cultural_sensitivity", 6),
    @("test_allergy_scenarios.py", "test_common_restriction_prompts", 16),
    @("test_allergy_scenarios.py", "test_peanut_allergy", 25),
    @("test_allergy_scenarios.py", "test_gluten_free", 29),
    @("test_allergy_scenarios.py", "test_dairy_free", 34),
    @("test_allergy_scenarios.py", "test_religious_restriction", 39),
    @("test_allergy_scenarios.py", "test_multi_axis_prompts", 52),
    @("test_allergy_scenarios.py", "test_ambiguous_cases", 67),
    @("test_equipment_scenarios.py", "test_common_equipment_prompts", 17),
    @("test_equipment_scenarios.py", "test_instant_pot_chicken_rice", 25),
    @("test_equipment_scenarios.py", "test_slow_cooker_soup", 29),
    @("test_equipment_scenarios.py", "test_oven_bake_conflict", 33),
    @("test_equipment_scenarios.py", "test_air_fryer_dessert", 38),
    @("test_equipment_scenarios.py", "test_wok_stirfry", 42),
    @("test_equipment_scenarios.py", "test_microwave_restriction", 47),
    @("test_response_logic.py", "test_known_ingredients", 4),
    @("test_response_logic.py", "test_equipment_detection", 8),
    @("test_response_logic.py", "test_restricted_ingredient_response", 12),
    @("test_response_logic.py", "test_unsure_user_input", 16),
    @("test_response_logic.py", "test_empty_input", 20),
    @("test_response_logic.py", "test_bias_trigger", 24)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $entry = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $entry[0]
    $ws.Cells.Item($rowNum, 2).Value = $entry[1]
    if ($rowNum -eq 1) {
        $ws.Cells.Item($rowNum, 3).Value = "Line Number"
    } else {
        $ws.Cells.Item($rowNum, 3).Value = $entry[2]
    }
}
